# Add new columns I (I0) and J (IF) to the sheet, mirroring the existing
# H (IP) column's header style, and fill in the data rows 2-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Copy the formatting of the existing "IP" header (H1) onto the two new
# header cells so they pick up the same bold/border/alignment style, then
# set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (rows 2-25) ---
# Values for column I (I0) and column J (IF) per row.
$data = @{
    2  = @(1, 6)
    3  = @(1, 5)
    4  = @(1, 6)
    5  = @(1, 5)
    6  = @(5, 7)
    7  = @(1, 3)
    8  = @(1, 6)
    9  = @(1, 4)
    10 = @(1, 6)
    11 = @(1, 5)
    12 = @(1, 6)
    13 = @(1, 4)
    14 = @(1, 5)
    15 = @(1, 3)
    16 = @(1, 6)
    17 = @(1, 6)
    18 = @(1, 6)
    19 = @(1, 5)
    20 = @(1, 5)
    21 = @(1, 7)
    22 = @(1, 6)
    23 = @(1, 5)
    24 = @(1, 5)
    25 = @(1, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
